# Generate Report for handoff
# Insert a new row for file "abb1d472-4017-4361-85ae-83523e266540" (status
# "Ready for handoff") right before the "ac6c1e5d-c464-49b4-8cda-02dce0384b70"
# row on all three worksheets (Overview, zh-cn, de-de), shifting the rows
# below it down by one, and fix up the hyperlinks to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Insert a new row 7 (pushes the old rows 7,8,9 down to 8,9,10); the new
# row inherits the style of the row above it (A=s1, B/C default).
$ws1.Rows.Item(7).Insert()

$ws1.Range("A7").Value = "abb1d472-4017-4361-85ae-83523e266540.md"
$ws1.Range("B7").Value = "Ready for handoff"
$ws1.Range("C7").Value = "Ready for handoff"

# Rebuild all hyperlinks on this sheet in the correct order/targets, since
# inserting a row does not shift the existing hyperlink anchors.
$ws1.Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d16bf8edabe688209a9175f0df403bff4b3ee8f5/e2e/bcf63ebd-162c-4609-9481-49dfaa24780a.md", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/760c6e0b3a33fd777a23776cb767086dff71dc1f/e2e/03acf302-3652-4579-973c-b89bbd18ffc4.md", "", "", "03acf302-3652-4579-973c-b89bbd18ffc4.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c174fd49aced88769ee5a2cb675486f334333fc9/e2e/09454cc8-cb47-49dc-8aac-94922237deaa.md", "", "", "09454cc8-cb47-49dc-8aac-94922237deaa.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c174fd49aced88769ee5a2cb675486f334333fc9/e2e/b776789b-d619-4ed8-8aac-5fd941d4d17c.md", "", "", "b776789b-d619-4ed8-8aac-5fd941d4d17c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/34f8effd44a64666bf79f84c4278f6e9de447d25/e2e/d3db4972-7764-4750-993a-1277e12b9ea9.md", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/cec6c615020baaab5166888f1fabfb9651101e33/e2e/abb1d472-4017-4361-85ae-83523e266540.md", "", "", "abb1d472-4017-4361-85ae-83523e266540.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/9149a40a9b28279b30d9bf2a9317e6b19baae428/e2e/ac6c1e5d-c464-49b4-8cda-02dce0384b70.md", "", "", "ac6c1e5d-c464-49b4-8cda-02dce0384b70.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/d6932189bf5e61fee79365e7e64c8782094aca2a/e2e/c565775a-50fe-4bf9-9527-4d2387950bf9.md", "", "", "c565775a-50fe-4bf9-9527-4d2387950bf9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/d16bf8edabe688209a9175f0df403bff4b3ee8f5/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(7).Insert()

$ws2.Range("A7").Value = "abb1d472-4017-4361-85ae-83523e266540.md"
$ws2.Range("B7").Value = "Ready for handoff"
$ws2.Range("C7").Value = "abb1d472-4017-4361-85ae-83523e266540.cec6c615020baaab5166888f1fabfb9651101e33.zh-cn.xlf"
$ws2.Range("D7").Value = "2016-01-25 05:45:32"
$ws2.Range("G7").Value = "0001-01-01 00:00:00"
$ws2.Range("H7").Value = "Include"

$ws2.Hyperlinks.Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d16bf8edabe688209a9175f0df403bff4b3ee8f5/e2e/bcf63ebd-162c-4609-9481-49dfaa24780a.md", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d283d063a838785fb0dfd735a91d3647231def72/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.zh-cn.xlf", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b8f602e688c481d1cde57dc33e1de70c82d8df5c/e2e/bcf63ebd-162c-4609-9481-49dfaa24780a.md", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/588af6261b647bc52d838c8322aa16a22e6fc01d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.zh-cn.xlf", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/760c6e0b3a33fd777a23776cb767086dff71dc1f/e2e/03acf302-3652-4579-973c-b89bbd18ffc4.md", "", "", "03acf302-3652-4579-973c-b89bbd18ffc4.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/471c0272f92579ae7195ad4a33ffb470193db388/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/03acf302-3652-4579-973c-b89bbd18ffc4.a1b011afc322fa2ea2e491f688a206adf34c2fc5.zh-cn.xlf", "", "", "03acf302-3652-4579-973c-b89bbd18ffc4.a1b011afc322fa2ea2e491f688a206adf34c2fc5.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c174fd49aced88769ee5a2cb675486f334333fc9/e2e/09454cc8-cb47-49dc-8aac-94922237deaa.md", "", "", "09454cc8-cb47-49dc-8aac-94922237deaa.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d25626040259023ca675ca378693f1f437e3b53a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/09454cc8-cb47-49dc-8aac-94922237deaa.4f308d36d7eafee632b9ea1f4adb2e8c37a526eb.zh-cn.xlf", "", "", "09454cc8-cb47-49dc-8aac-94922237deaa.4f308d36d7eafee632b9ea1f4adb2e8c37a526eb.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c174fd49aced88769ee5a2cb675486f334333fc9/e2e/b776789b-d619-4ed8-8aac-5fd941d4d17c.md", "", "", "b776789b-d619-4ed8-8aac-5fd941d4d17c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d25626040259023ca675ca378693f1f437e3b53a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/b776789b-d619-4ed8-8aac-5fd941d4d17c.dbfc76c3a5bc4c0376d60c05468c0ddcb5cd4b4e.zh-cn.xlf", "", "", "b776789b-d619-4ed8-8aac-5fd941d4d17c.dbfc76c3a5bc4c0376d60c05468c0ddcb5cd4b4e.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/34f8effd44a64666bf79f84c4278f6e9de447d25/e2e/d3db4972-7764-4750-993a-1277e12b9ea9.md", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9a9c985454eb30bdceb85ec575239fd6c7964c47/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.zh-cn.xlf", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2fcccbea2fa1d7bec1ed8c85a583ef39c79a9938/e2e/d3db4972-7764-4750-993a-1277e12b9ea9.md", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2a0332230824ef29b14e7b1ef75155ad87d64624/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.zh-cn.xlf", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/cec6c615020baaab5166888f1fabfb9651101e33/e2e/abb1d472-4017-4361-85ae-83523e266540.md", "", "", "abb1d472-4017-4361-85ae-83523e266540.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cec6c615020baaab5166888f1fabfb9651101e33/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/abb1d472-4017-4361-85ae-83523e266540.cec6c615020baaab5166888f1fabfb9651101e33.zh-cn.xlf", "", "", "abb1d472-4017-4361-85ae-83523e266540.cec6c615020baaab5166888f1fabfb9651101e33.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/9149a40a9b28279b30d9bf2a9317e6b19baae428/e2e/ac6c1e5d-c464-49b4-8cda-02dce0384b70.md", "", "", "ac6c1e5d-c464-49b4-8cda-02dce0384b70.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1a4f8e56516b74650c9aaca497b54ad17ff1ab0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ac6c1e5d-c464-49b4-8cda-02dce0384b70.b44c68e17ae6ad15d3d3f0f0e344be729f4e1235.zh-cn.xlf", "", "", "ac6c1e5d-c464-49b4-8cda-02dce0384b70.b44c68e17ae6ad15d3d3f0f0e344be729f4e1235.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/d6932189bf5e61fee79365e7e64c8782094aca2a/e2e/c565775a-50fe-4bf9-9527-4d2387950bf9.md", "", "", "c565775a-50fe-4bf9-9527-4d2387950bf9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5dc647c53c849ea0f82a431b8bf4d638b2d81eb0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/c565775a-50fe-4bf9-9527-4d2387950bf9.b7bf726f5cdc4dfd77ee110198b844ddd68c40fd.zh-cn.xlf", "", "", "c565775a-50fe-4bf9-9527-4d2387950bf9.b7bf726f5cdc4dfd77ee110198b844ddd68c40fd.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/d16bf8edabe688209a9175f0df403bff4b3ee8f5/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(7).Insert()

$ws3.Range("A7").Value = "abb1d472-4017-4361-85ae-83523e266540.md"
$ws3.Range("B7").Value = "Ready for handoff"
$ws3.Range("C7").Value = "abb1d472-4017-4361-85ae-83523e266540.cec6c615020baaab5166888f1fabfb9651101e33.de-de.xlf"
$ws3.Range("D7").Value = "2016-01-25 05:45:43"
$ws3.Range("G7").Value = "0001-01-01 00:00:00"
$ws3.Range("H7").Value = "Include"

$ws3.Hyperlinks.Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d16bf8edabe688209a9175f0df403bff4b3ee8f5/e2e/bcf63ebd-162c-4609-9481-49dfaa24780a.md", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f936bfff0148d53da167d095f2912a503bf0bbc3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.de-de.xlf", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/5e342f9a30e064f0886a79178eaf0f287ce923ed/e2e/bcf63ebd-162c-4609-9481-49dfaa24780a.md", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ab3f6f3b4e92fdee92145b0aa7c1261ea432712b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.de-de.xlf", "", "", "bcf63ebd-162c-4609-9481-49dfaa24780a.97de180362651a6615f1cbfaa76e33ca0d961129.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/760c6e0b3a33fd777a23776cb767086dff71dc1f/e2e/03acf302-3652-4579-973c-b89bbd18ffc4.md", "", "", "03acf302-3652-4579-973c-b89bbd18ffc4.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/01d9973b672f954b06c8a920b05a9fff8f76297e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/03acf302-3652-4579-973c-b89bbd18ffc4.a1b011afc322fa2ea2e491f688a206adf34c2fc5.de-de.xlf", "", "", "03acf302-3652-4579-973c-b89bbd18ffc4.a1b011afc322fa2ea2e491f688a206adf34c2fc5.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c174fd49aced88769ee5a2cb675486f334333fc9/e2e/09454cc8-cb47-49dc-8aac-94922237deaa.md", "", "", "09454cc8-cb47-49dc-8aac-94922237deaa.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/111d40962d8a64ea9b1426333d6970210bd7e572/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/09454cc8-cb47-49dc-8aac-94922237deaa.4f308d36d7eafee632b9ea1f4adb2e8c37a526eb.de-de.xlf", "", "", "09454cc8-cb47-49dc-8aac-94922237deaa.4f308d36d7eafee632b9ea1f4adb2e8c37a526eb.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/c174fd49aced88769ee5a2cb675486f334333fc9/e2e/b776789b-d619-4ed8-8aac-5fd941d4d17c.md", "", "", "b776789b-d619-4ed8-8aac-5fd941d4d17c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/111d40962d8a64ea9b1426333d6970210bd7e572/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/b776789b-d619-4ed8-8aac-5fd941d4d17c.dbfc76c3a5bc4c0376d60c05468c0ddcb5cd4b4e.de-de.xlf", "", "", "b776789b-d619-4ed8-8aac-5fd941d4d17c.dbfc76c3a5bc4c0376d60c05468c0ddcb5cd4b4e.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/34f8effd44a64666bf79f84c4278f6e9de447d25/e2e/d3db4972-7764-4750-993a-1277e12b9ea9.md", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed8bcc4523e6cd2272436a8ddafc22e10175a30a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.de-de.xlf", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/95e611855d0c45b94b1c23792623481e007dd842/e2e/d3db4972-7764-4750-993a-1277e12b9ea9.md", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/38912e11328abc9336b33cc4e5a9a794cfd741d3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.de-de.xlf", "", "", "d3db4972-7764-4750-993a-1277e12b9ea9.a036ecc265d98c35051ae0e03f165862d230ab0d.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/cec6c615020baaab5166888f1fabfb9651101e33/e2e/abb1d472-4017-4361-85ae-83523e266540.md", "", "", "abb1d472-4017-4361-85ae-83523e266540.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cec6c615020baaab5166888f1fabfb9651101e33/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/abb1d472-4017-4361-85ae-83523e266540.cec6c615020baaab5166888f1fabfb9651101e33.de-de.xlf", "", "", "abb1d472-4017-4361-85ae-83523e266540.cec6c615020baaab5166888f1fabfb9651101e33.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/9149a40a9b28279b30d9bf2a9317e6b19baae428/e2e/ac6c1e5d-c464-49b4-8cda-02dce0384b70.md", "", "", "ac6c1e5d-c464-49b4-8cda-02dce0384b70.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90ca176103aad7cbb5665acdadb446a0fd8a0dd5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ac6c1e5d-c464-49b4-8cda-02dce0384b70.b44c68e17ae6ad15d3d3f0f0e344be729f4e1235.de-de.xlf", "", "", "ac6c1e5d-c464-49b4-8cda-02dce0384b70.b44c68e17ae6ad15d3d3f0f0e344be729f4e1235.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/d6932189bf5e61fee79365e7e64c8782094aca2a/e2e/c565775a-50fe-4bf9-9527-4d2387950bf9.md", "", "", "c565775a-50fe-4bf9-9527-4d2387950bf9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e35d64dd1a2f3e71ba4f18686bb76fd93cb0ea08/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/c565775a-50fe-4bf9-9527-4d2387950bf9.b7bf726f5cdc4dfd77ee110198b844ddd68c40fd.de-de.xlf", "", "", "c565775a-50fe-4bf9-9527-4d2387950bf9.b7bf726f5cdc4dfd77ee110198b844ddd68c40fd.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A10"), "https://github.com/OpenLocalizationTest/oltest/blob/d16bf8edabe688209a9175f0df403bff4b3ee8f5/.localization-config", "", "", ".localization-config") | Out-Null
